$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.595.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.478.38"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9616"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "280.80"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3667"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3082"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40.07"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.061"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06669"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.007"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.527"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.09"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.220"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9641"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001033"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.478.22"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05978"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.48"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.643.89"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.97"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.117"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -8.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.29"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.640.13"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.98"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.978"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.047"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8183"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07997"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.43%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.219"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05810"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.737"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02049"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9632"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.41"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.527"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1880"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5317"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.33"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.546"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.54"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5209"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.829"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06504"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9915"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.43%  "
